$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels from Russian to English
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = "Income"

# Update the active selection (was A4:G17 with active cell G4, now just A3)
$ws.Activate()
$ws.Range("A3").Select()
